$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was updated from
# 45186 (2023-09-17) to 45188 (2023-09-19) for every data row (2 through 295).
$ws.Range("C2:C295").Value = 45188
